# Updates cryptos list prices and volume percentages (and restores the
# FraxShare/VeChain row order) to match the latest coinranking.com snapshot,
# mirroring the automated GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "40.074.21"
$ws.Range("E2").Value = "  +2.33%  "
$ws.Range("D3").Value = "2.234.83"
$ws.Range("E3").Value = "  +0.93%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'294.50"
$ws.Range("E5").Value = "  -0.70%  "
$ws.Range("D6").Value = "'86.94"
$ws.Range("E6").Value = "  +5.24%  "
$ws.Range("E7").Value = "  +1.34%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").Value = "'0.473"
$ws.Range("E9").Value = "  +1.52%  "
$ws.Range("D10").Value = "'31.00"
$ws.Range("E10").Value = "  +7.29%  "
$ws.Range("D11").Value = "'0.0791"
$ws.Range("E11").Value = "  +2.37%  "
$ws.Range("D12").Value = "'46.98"
$ws.Range("E12").Value = "  -0.08%  "
$ws.Range("E13").Value = "  +1.32%  "
$ws.Range("E14").Value = "  +3.27%  "
$ws.Range("D15").Value = "2.583.75"
$ws.Range("E15").Value = "  +0.91%  "
$ws.Range("D16").Value = "'14.12"
$ws.Range("E16").Value = "  +0.29%  "
$ws.Range("D17").Value = "2.249.71"
$ws.Range("E17").Value = "  +1.22%  "
$ws.Range("E18").Value = "  +2.64%  "
$ws.Range("D19").Value = "39.998.29"
$ws.Range("E19").Value = "  +2.35%  "
$ws.Range("E20").Value = "  +2.66%  "
$ws.Range("D21").Value = "'11.25"
$ws.Range("E21").Value = "  +10.93%  "
$ws.Range("E22").Value = "  +2.23%  "
$ws.Range("D23").Value = "'65.48"
$ws.Range("E23").Value = "  +1.15%  "
$ws.Range("D24").Value = "'235.66"
$ws.Range("E24").Value = "  +3.87%  "
$ws.Range("E25").Value = "  +0.04%  "
$ws.Range("E26").Value = "  +3.10%  "
$ws.Range("D27").Value = "'1.85"
$ws.Range("E27").Value = "  +4.88%  "
$ws.Range("D28").Value = "'22.84"
$ws.Range("E28").Value = "  +1.50%  "
$ws.Range("E29").Value = "  +2.69%  "
$ws.Range("E30").Value = "  +3.11%  "
$ws.Range("D31").Value = "'33.30"
$ws.Range("E31").Value = "  +4.82%  "
$ws.Range("D32").Value = "'152.88"
$ws.Range("E32").Value = "  +2.77%  "
$ws.Range("D33").Value = "'0.999"
$ws.Range("E33").Value = "  -0.22%  "
$ws.Range("D34").Value = "'4.94"
$ws.Range("E34").Value = "  +2.56%  "
$ws.Range("D35").Value = "'0.0722"
$ws.Range("E35").Value = "  +3.84%  "
$ws.Range("E36").Value = "  +3.13%  "
$ws.Range("D37").Value = "'16.27"
$ws.Range("E37").Value = "  +9.95%  "
$ws.Range("E38").Value = "  +5.80%  "
$ws.Range("E39").Value = "  +5.18%  "
$ws.Range("D40").Value = "'0.111"
$ws.Range("E40").Value = "  +2.36%  "
$ws.Range("D41").Value = "'1.71"
$ws.Range("E41").Value = "  +5.89%  "
$ws.Range("D42").Value = "'3.82"
$ws.Range("E42").Value = "  +3.82%  "
$ws.Range("D43").Value = "2.050.05"
$ws.Range("E43").Value = "  +7.32%  "
$ws.Range("E44").Value = "  +6.76%  "
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").Value = "'0.0271"
$ws.Range("E45").Value = "  +5.32%  "
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").Value = "'9.98"
$ws.Range("E46").Value = "  +12.26%  "
$ws.Range("D47").Value = "'17.09"
$ws.Range("E47").Value = "  +7.86%  "
$ws.Range("D48").Value = "'2.61"
$ws.Range("E48").Value = "  -0.83%  "
$ws.Range("D49").Value = "2.444.76"
$ws.Range("E49").Value = "  +0.59%  "
$ws.Range("D50").Value = "'71.48"
$ws.Range("E50").Value = "  +2.11%  "
$ws.Range("D51").Value = "'89.36"
$ws.Range("E51").Value = "  +2.89%  "
